# Update column F (dSF) values for rows 2-20 (except rows 13 and 15, which
# are unchanged) to reflect a repull of the data / recalculated mean.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = 0
    4  = 2
    5  = 4
    6  = 1
    7  = 6
    8  = -3
    9  = -1
    10 = 4
    11 = 9
    12 = -4
    14 = -2
    16 = -1
    17 = 9
    18 = 7
    19 = -1
    20 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
